# Append the next day's gold-price row to the bottom of the data table
# (mirrors the daily GitHub Actions scrape-and-commit job).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$dateCell = $ws.Cells.Item($newRow, 1)
$priceCell = $ws.Cells.Item($newRow, 2)

# Write the new date through a literal-text formula first so Excel's
# input-parsing doesn't silently reinterpret the dd-mm-yyyy-looking
# string as a date serial; then collapse the formula to its plain
# string result in place so formatting/style stay untouched.
$dateCell.Formula = "=""07-12-2025"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

# Today's price paragraph is unchanged from yesterday's on the source
# site, so copy the prior row's text verbatim.
$priceText = $ws.Cells.Item($lastRow, 2).Value2
$priceCell.Value = $priceText

$excel.CutCopyMode = $false
